$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("SchemaOrganization")
$ws.Range("B2").Value = "http://example.com/organization1:Image1"

$ws = $wb.Worksheets.Item("CreativeCommonsLicense")
$ws.Range("A4").Value = "http://creativecommons.org/licenses/by-sa/2.0/"
$ws.Range("G4").Value = "by-sa"
$ws.Range("H4").Value = "Attribution-ShareAlike 2.0 Generic"
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "2.0"
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = ""
$ws.Range("L4").Value = ""
$ws.Range("A5").Value = "http://creativecommons.org/licenses/nc/1.0/"
$ws.Range("G5").Value = ""
$ws.Range("H5").Value = ""
$ws.Range("I5").Value = "nc"
$ws.Range("J5").Value = "NonCommercial 1.0 Generic"
$ws.Range("K5").NumberFormat = "@"
$ws.Range("K5").Value = "1.0"

$ws = $wb.Worksheets.Item("RightsStatementsDotOrgRightsStatement")
$ws.Range("E2").Value = "You may need to obtain other permissions for your intended use. For example, other rights such as publicity, privacy or moral rights may limit how you may use the material."
$ws.Range("E3").Value = "You may need to obtain other permissions for your intended use. For example, other rights such as publicity, privacy or moral rights may limit how you may use the material."

$ws = $wb.Worksheets.Item("FoafPerson")
$ws.Range("E4").Value = "http://example.com/person4:Image0"

$ws = $wb.Worksheets.Item("SchemaPerson")
$ws.Range("D2").Value = "http://example.com/person1:Image1"
$ws.Range("D3").Value = "http://example.com/person3:Image1"

$ws = $wb.Worksheets.Item("RdfProperty")
$ws.Range("C3").Value = "dcterms:extent:Image1"
$ws.Range("C5").Value = "dcterms:medium:Image0"
$ws.Range("C6").Value = "dcterms:publisher:Image1"
$ws.Range("C7").Value = "dcterms:source:Image1"
$ws.Range("C10").Value = "dcterms:title:Image1"

$ws = $wb.Worksheets.Item("SchemaProperty")
$ws.Range("C3").Value = "schema:name:Image1"
$ws.Range("C4").Value = "schema:spatial:Image1"

$ws = $wb.Worksheets.Item("CmsPropertyGroup")
$ws.Range("C2").Value = "urn:paradicms_etl:pipeline:synthetic_data:property_group:Image1"

$ws = $wb.Worksheets.Item("SkosConcept")
$ws.Range("B3").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:1:Image1"
$ws.Range("B6").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:4:Image1"
$ws.Range("B7").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:5:Image0"
$ws.Range("B8").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:6:Image0"
$ws.Range("B10").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:8:Image1"
$ws.Range("B11").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:9:Image0"
$ws.Range("B15").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:13:Image0"
$ws.Range("B16").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:14:Image1"
$ws.Range("B17").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:15:Image0"
$ws.Range("B18").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:16:Image1"
$ws.Range("B20").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:18:Image0"
$ws.Range("B23").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:21:Image1"
$ws.Range("B24").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:22:Image0"
$ws.Range("B26").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:24:Image1"
$ws.Range("B27").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:25:Image1"
$ws.Range("B29").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:27:Image0"
$ws.Range("B30").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:28:Image1"
$ws.Range("B33").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:31:Image0"
$ws.Range("B37").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:35:Image1"
$ws.Range("B38").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:36:Image1"
$ws.Range("B40").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:38:Image1"
$ws.Range("B41").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:39:Image1"
$ws.Range("B43").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:41:Image0"
$ws.Range("B45").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:43:Image0"
$ws.Range("B48").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:46:Image0"
$ws.Range("B51").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:49:Image0"
$ws.Range("B52").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:50:Image1"
$ws.Range("B53").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:51:Image1"
$ws.Range("B55").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:53:Image0"
$ws.Range("B57").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:55:Image0"
$ws.Range("B62").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:60:Image0"
$ws.Range("B66").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:64:Image0"
$ws.Range("B67").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:65:Image1"
$ws.Range("B68").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:66:Image1"
$ws.Range("B69").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:67:Image1"
$ws.Range("B71").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:69:Image1"
$ws.Range("B73").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:71:Image1"
$ws.Range("B74").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:72:Image0"
$ws.Range("B75").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:73:Image1"
$ws.Range("B76").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:74:Image1"
$ws.Range("B78").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:76:Image0"
$ws.Range("B80").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:78:Image1"

$ws = $wb.Worksheets.Item("SchemaDefinedTerm")
$ws.Range("B2").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:80:Image1"
$ws.Range("B4").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:82:Image0"
$ws.Range("B8").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:86:Image0"
$ws.Range("B9").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:87:Image1"

$ws = $wb.Worksheets.Item("FoafOrganization")
$ws.Range("C4").Value = "http://example.com/organization4:Image0"
